# Implement error message display - appends a new row (row 79) of sensor
# data to each of the four worksheets, matching the latest reading that
# had previously been filtered/trimmed from the dataset.

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = "2025-03-07 14:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = "2025-03-07 14:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = "2025-03-07 14:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = "2025-03-07 14:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F

    # Column G holds a 24-digit identifier that exceeds Excel's numeric
    # precision, so it must be stored as text (as the rest of column G is)
    # rather than being auto-converted to a floating point number.
    $gCell = $ws.Cells.Item($newRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $rowData.G

    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}

Write-Host "Added row 79 to all four sheets"
